# Updated cryptos list on Wed Apr  3 19:25:51 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.844.55"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "3.317.48"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "557.91"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").Value = "185.82"
$ws.Range("E6").Value = "  +0.58%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").Value = "3.310.56"
$ws.Range("E8").Value = "  +1.40%  "

$ws.Range("E9").Value = "  -2.52%  "

$ws.Range("E10").Value = "  -6.70%  "

$ws.Range("E11").Value = "  -1.48%  "

$ws.Range("D12").Value = "45.84"
$ws.Range("E12").Value = "  -3.10%  "

$ws.Range("E13").Value = "  -1.50%  "

$ws.Range("D14").Value = "3.850.55"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("D16").Value = "571.77"
$ws.Range("E16").Value = "  -9.20%  "

$ws.Range("D17").Value = "65.850.66"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").Value = "3.328.69"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").Value = "17.65"
$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("D21").Value = "10.84"
$ws.Range("E21").Value = "  -4.37%  "

$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("D23").Value = "'18.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  +1.42%  "

$ws.Range("D25").Value = "97.73"
$ws.Range("E25").Value = "  -7.95%  "

$ws.Range("D26").Value = "3.94"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("D30").Value = "30.49"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").Value = "6.71"
$ws.Range("E31").Value = "  +7.54%  "

# Rows 32/33 swap: Bittensor <-> dogwifhat with new values
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "3.68"
$ws.Range("E32").Value = "  -8.65%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "561.64"
$ws.Range("E33").Value = "  +4.27%  "

$ws.Range("D34").Value = "10.82"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("E35").Value = "  -1.70%  "

$ws.Range("D36").Value = "3.741.95"
$ws.Range("E36").Value = "  +1.03%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "55.49"
$ws.Range("E38").Value = "  -3.43%  "

$ws.Range("D39").Value = "'33.80"
$ws.Range("E39").Value = "  +3.28%  "

$ws.Range("E40").Value = "  -3.83%  "

$ws.Range("D41").Value = "0.0$([char]0x2083)0686"
$ws.Range("E41").Value = "  -5.84%  "

$ws.Range("E42").Value = "  -5.54%  "

$ws.Range("E43").Value = "  -8.59%  "

$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("D47").Value = "2.99"
$ws.Range("E47").Value = "  -12.38%  "

# Rows 48/49 swap: FirstDigitalUSD <-> Stellar with new values
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.126"
$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  -3.35%  "

$ws.Range("D51").Value = "124.31"
$ws.Range("E51").Value = "  +1.70%  "
